$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)

# 1) Move the "Thanks for listening !" title shape up.
$title = $s.Shapes.Item(1)
$title.Top = 111.34055118110236

# 2) Add the new Github-link textbox below the title.
$tb = $s.Shapes.AddTextbox(1, 118.62992125984252, 197.48425196850394, 482.74015748031496, 93.4251968503937)
$tb.Name = "Google Shape;189;p21"

# No fill / no line, like the other shapes on this slide.
$tb.Fill.Visible = 0
$tb.Line.Visible = 0

# Body text properties (match the Google Slides export defaults used elsewhere
# in this deck: 91425 EMU insets, top anchored, no autofit, word wrap on).
$tf = $tb.TextFrame
$tf.WordWrap = -1
$tf.AutoSize = 0
$tf.MarginLeft = 7.198818897637795
$tf.MarginRight = 7.198818897637795
$tf.MarginTop = 7.198818897637795
$tf.MarginBottom = 7.198818897637795
$tf.VerticalAnchor = 1
$tf.HorizontalAnchor = 0

$tr = $tf.TextRange
$tr.Text = "https://github.com/PabloBerenguel/CAu-Natural-language-processing"
$tr.ParagraphFormat.Alignment = 2
$tr.ParagraphFormat.Bullet.Visible = 0
$tr.ParagraphFormat.SpaceBefore = 0
$tr.ParagraphFormat.SpaceAfter = 0
$tr.IndentLevel = 1
$tr.Font.LanguageID = "fr"
$tr.Font.Color.RGB = 16777215
